# Updating filtered feeds from workflow
# Appends a new feed entry as the next row of the "Filtered Feeds" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 87

$link = "https://www.fiercebiotech.com/medtech/guardant-health-gets-fda-nod-expand-use-blood-test-colon-cancer"
$keywords = "CDx"
$title = '<a href="https://www.fiercebiotech.com/medtech/guardant-health-gets-fda-nod-expand-use-blood-test-colon-cancer" hreflang="en">Guardant Health gets FDA nod to expand use of blood test for colon cancer</a>'

# Column A (link) is a hyperlink, same as every other row in the sheet.
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $link, [Type]::Missing, [Type]::Missing, [Type]::Missing)
$ws.Cells.Item($newRow, 1).Style = "Hyperlink"

# Column B (keywords) and column C (title) are plain text values.
$ws.Cells.Item($newRow, 2).Value2 = $keywords
$ws.Cells.Item($newRow, 3).Value2 = $title

Write-Host "Added row $newRow for $link"
